$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; Excel shifts all existing rows (33..116) down to (34..117)
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with its data
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44965
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 37000
$ws.Range("L33").Value = 38000
$ws.Range("M33").Value = 37500
$ws.Range("N33").Value = "$/malla 25 kilos"
$ws.Range("O33").Value = "Provincia del Elquí"
$ws.Range("P33").Value = 1500
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
